# Insert a new weekly price record for "Macroferia Regional de Talca - Cilantro".
# This pushes the existing rows 108..125 down to 109..126 and fills the newly
# opened row 108 with the latest observation (matches the commit's weekly
# fruit/vegetable data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 108 (shifts 108..125 -> 109..126).
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record.
$ws.Cells.Item(108, 1).Value  = 5
$ws.Cells.Item(108, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(108, 3).Value  = "Maule"
$ws.Cells.Item(108, 4).Value  = 45180
$ws.Cells.Item(108, 5).Value  = 7
$ws.Cells.Item(108, 6).Value  = 100112040
$ws.Cells.Item(108, 7).Value  = "Cilantro"
$ws.Cells.Item(108, 8).Value  = "Sin especificar"
$ws.Cells.Item(108, 9).Value  = "Primera"
$ws.Cells.Item(108, 10).Value = 150
$ws.Cells.Item(108, 11).Value = 9000
$ws.Cells.Item(108, 12).Value = 9000
$ws.Cells.Item(108, 13).Value = 9000
$ws.Cells.Item(108, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(108, 15).Value = "Región Metropolitana"
$ws.Cells.Item(108, 16).Value = 250
$ws.Cells.Item(108, 17).Value = 36
$ws.Cells.Item(108, 18).Value = "Hortaliza"
